# Update 08/10/2024
# The sheet contained a "link" header in A1 followed by 52 Facebook post
# URLs in A2:A53 (pulled from sharedStrings). This batch of links is removed
# from the list - delete those 52 rows entirely so everything below shifts
# up, leaving just the "link" header and the already-blank/formatted rows
# that followed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 2 through 53 (the 52 link entries), shifting remaining rows up.
$ws.Rows("2:53").Delete()

# Reflect the last on-screen selection/cursor position after the edit.
$ws.Range("I11").Select() | Out-Null
